$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numbers formatted as plain text (e.g. "1.003", "30.330.55"),
# so force text storage first -- otherwise Excel auto-converts numeric-looking strings into
# real numbers (and mangles multi-dot values / drops formatting), which is not what the
# source data feed produces.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "30.330.55"
$ws.Range("E2").Value = "  -2.06%  "

# Row 3
$ws.Range("D3").Value = "1.901.74"
$ws.Range("E3").Value = "  -2.81%  "

# Row 4
$ws.Range("D4").Value = "1.003"

# Row 5
$ws.Range("D5").Value = "237.46"
$ws.Range("E5").Value = "  -2.69%  "

# Row 6
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.09%  "

# Row 7
$ws.Range("D7").Value = "0.4713"
$ws.Range("E7").Value = "  -2.81%  "

# Row 8
$ws.Range("D8").Value = "0.2813"
$ws.Range("E8").Value = "  -4.29%  "

# Row 9
$ws.Range("D9").Value = "0.06589"
$ws.Range("E9").Value = "  -6.22%  "

# Row 10
$ws.Range("D10").Value = "18.52"
$ws.Range("E10").Value = "  -5.91%  "

# Row 11
$ws.Range("D11").Value = "99.36"
$ws.Range("E11").Value = "  -7.38%  "

# Row 12
$ws.Range("D12").Value = "0.07709"
$ws.Range("E12").Value = "  -1.26%  "

# Row 13
$ws.Range("D13").Value = "1.914.60"
$ws.Range("E13").Value = "  -2.01%  "

# Row 14
$ws.Range("D14").Value = "5.135"
$ws.Range("E14").Value = "  -5.84%  "

# Row 15
$ws.Range("D15").Value = "0.6623"
$ws.Range("E15").Value = "  -5.28%  "

# Row 16
$ws.Range("D16").Value = "30.363.79"
$ws.Range("E16").Value = "  -2.03%  "

# Row 17
$ws.Range("D17").Value = "250.70"
$ws.Range("E17").Value = "  -10.65%  "

# Row 18
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  +0.13%  "

# Row 19
$ws.Range("D19").Value = "0.000007387"
$ws.Range("E19").Value = "  -5.26%  "

# Row 20
$ws.Range("D20").Value = "12.54"
$ws.Range("E20").Value = "  -5.67%  "

# Row 21
$ws.Range("D21").Value = "5.345"
$ws.Range("E21").Value = "  -3.73%  "

# Row 22
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  -0.02%  "

# Row 23
$ws.Range("D23").Value = "6.245"
$ws.Range("E23").Value = "  -3.85%  "

# Row 24
$ws.Range("D24").Value = "9.290"
$ws.Range("E24").Value = "  -5.37%  "

# Row 25
$ws.Range("D25").Value = "163.95"
$ws.Range("E25").Value = "  -2.81%  "

# Row 26
$ws.Range("D26").Value = "18.69"
$ws.Range("E26").Value = "  -5.63%  "

# Row 27
$ws.Range("D27").Value = "2.030"
$ws.Range("E27").Value = "  -6.60%  "

# Row 28
$ws.Range("D28").Value = "0.1005"
$ws.Range("E28").Value = "  -4.04%  "

# Row 29
$ws.Range("E29").Value = "  -0.68%  "

# Row 30
$ws.Range("D30").Value = "4.600"
$ws.Range("E30").Value = "  -0.12%  "

# Row 31
$ws.Range("D31").Value = "1.505"
$ws.Range("E31").Value = "  -4.15%  "

# Row 32
$ws.Range("D32").Value = "4.200"
$ws.Range("E32").Value = "  -5.32%  "

# Row 33
$ws.Range("D33").Value = "0.04694"
$ws.Range("E33").Value = "  -3.94%  "

# Row 34
$ws.Range("D34").Value = "0.7198"
$ws.Range("E34").Value = "  -3.49%  "

# Row 35
$ws.Range("D35").Value = "1.099"
$ws.Range("E35").Value = "  -5.71%  "

# Row 36
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  -0.68%  "

# Row 38
$ws.Range("D38").Value = "0.01896"
$ws.Range("E38").Value = "  -5.25%  "

# Row 39
$ws.Range("D39").Value = "2.595"
$ws.Range("E39").Value = "  -3.37%  "

# Row 40
$ws.Range("D40").Value = "6.197"
$ws.Range("E40").Value = "  -4.88%  "

# Row 41
$ws.Range("D41").Value = "72.15"
$ws.Range("E41").Value = "  -7.05%  "

# Row 42
$ws.Range("D42").Value = "1.963"
$ws.Range("E42").Value = "  -7.45%  "

# Row 43
$ws.Range("D43").Value = "106.02"
$ws.Range("E43").Value = "  -2.79%  "

# Row 44
$ws.Range("D44").Value = "0.8551"
$ws.Range("E44").Value = "  -5.16%  "

# Row 45
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.04%  "

# Row 46
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.034.33"
$ws.Range("E46").Value = "  +4.42%  "

# Row 47
$ws.Range("D47").Value = "0.4189"
$ws.Range("E47").Value = "  -5.71%  "

# Row 48
$ws.Range("D48").Value = "7.356"
$ws.Range("E48").Value = "  -8.20%  "

# Row 49
$ws.Range("D49").Value = "0.1184"
$ws.Range("E49").Value = "  -5.00%  "

# Row 50
$ws.Range("D50").Value = "34.31"
$ws.Range("E50").Value = "  -4.51%  "

# Row 51
$ws.Range("D51").Value = "8.774"
$ws.Range("E51").Value = "  -5.69%  "
